$wb = $excel.ActiveWorkbook
$origActiveSheet = $wb.ActiveSheet.Name
$ws = $wb.Worksheets.Item("Global")

# Add the new "FullName" column header next to the existing datasheet columns.
$ws.Range("G1").Value = "FullName"

# Give the new column G2 cell the same "closing" border style that F2 (the
# previous last column) currently has, then relax F2 back to an interior
# border now that it is no longer the last column.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("F2").Borders(10).LineStyle = -4142

# Size the new column and move the selection past it, like a value was
# just entered there and the cursor advanced to the next cell.
$ws.Columns("G").ColumnWidth = 8.61
[void]$ws.Range("H2").Select()

# Restore whichever sheet/tab was active before this edit.
[void]$wb.Worksheets.Item($origActiveSheet).Activate()
